$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A61 holds a date-looking string ("2025/10/05") that must stay plain text
# (matching the rest of column A), not get auto-converted into a date
# serial number by Excel's input parser. Force text storage via a
# temporary text number format, then drop the format again so the cell
# doesn't pick up a lingering style index.
$ws.Range("A61").NumberFormat = "@"
$ws.Range("A61").Value = "2025/10/05"
$ws.Range("A61").ClearFormats()

$ws.Range("B61").Value = "日"
$ws.Range("C61").Value = 0
$ws.Range("D61").Value = 46
